# Update the division-practice worksheet table: each data row (1, 5, 9,
# 13, 17) holds 5 "a÷b=" problems in columns 1-5. Replace the text of each
# cell directly (by position) rather than a global Find/Replace, since a
# few old values coincide with other cells' new values (e.g. "51÷9=" and
# "58÷2=" are both an old value in one cell and the new value in another),
# which would make a sequential whole-document replace ambiguous/unsafe.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "95÷7="
$t.Cell(1, 2).Range.Text  = "36÷8="
$t.Cell(1, 3).Range.Text  = "51÷9="
$t.Cell(1, 4).Range.Text  = "58÷2="
$t.Cell(1, 5).Range.Text  = "53÷9="

$t.Cell(5, 1).Range.Text  = "53÷6="
$t.Cell(5, 2).Range.Text  = "50÷2="
$t.Cell(5, 3).Range.Text  = "13÷9="
$t.Cell(5, 4).Range.Text  = "20÷6="
$t.Cell(5, 5).Range.Text  = "55÷5="

$t.Cell(9, 1).Range.Text  = "57÷8="
$t.Cell(9, 2).Range.Text  = "92÷5="
$t.Cell(9, 3).Range.Text  = "86÷3="
$t.Cell(9, 4).Range.Text  = "78÷4="
$t.Cell(9, 5).Range.Text  = "23÷2="

$t.Cell(13, 1).Range.Text = "82÷7="
$t.Cell(13, 2).Range.Text = "23÷6="
$t.Cell(13, 3).Range.Text = "43÷2="
$t.Cell(13, 4).Range.Text = "77÷2="
$t.Cell(13, 5).Range.Text = "95÷3="

$t.Cell(17, 1).Range.Text = "24÷9="
$t.Cell(17, 2).Range.Text = "33÷5="
$t.Cell(17, 3).Range.Text = "39÷6="
$t.Cell(17, 4).Range.Text = "96÷3="
$t.Cell(17, 5).Range.Text = "53÷9="

Write-Output "Updated 25 division problems."
